$d = $word.ActiveDocument

# Update the date heading (first paragraph)
$d.Content.Find.Execute("2025-02-17 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-02-18 Tuesday", 2) | Out-Null

# Update the division problems in the table. Cell addressing (1-based
# row/column) is used instead of text search-and-replace because one of
# the new values ("186÷4=") collides with an old value elsewhere in the
# table, which would make a blind global replace ambiguous/order-dependent.
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text  = "990÷8="
$t.Cell(1,2).Range.Text  = "214÷3="
$t.Cell(1,3).Range.Text  = "368÷9="
$t.Cell(1,4).Range.Text  = "512÷2="
$t.Cell(1,5).Range.Text  = "708÷2="

$t.Cell(5,1).Range.Text  = "427÷4="
$t.Cell(5,2).Range.Text  = "397÷4="
$t.Cell(5,3).Range.Text  = "840÷3="
$t.Cell(5,4).Range.Text  = "511÷7="
$t.Cell(5,5).Range.Text  = "240÷6="

$t.Cell(9,1).Range.Text  = "372÷7="
$t.Cell(9,2).Range.Text  = "548÷8="
$t.Cell(9,3).Range.Text  = "538÷7="
$t.Cell(9,4).Range.Text  = "810÷8="
$t.Cell(9,5).Range.Text  = "224÷5="

$t.Cell(13,1).Range.Text = "366÷9="
$t.Cell(13,2).Range.Text = "236÷9="
$t.Cell(13,3).Range.Text = "847÷2="
$t.Cell(13,4).Range.Text = "194÷2="
$t.Cell(13,5).Range.Text = "439÷4="

$t.Cell(17,1).Range.Text = "776÷7="
$t.Cell(17,2).Range.Text = "562÷4="
$t.Cell(17,3).Range.Text = "318÷7="
$t.Cell(17,4).Range.Text = "186÷4="
$t.Cell(17,5).Range.Text = "763÷5="

Write-Host "Done"
